# Chapter_4_Table_S4.17.xlsx edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header renames (row 2) - also syncs the Table17 column names since
#    the header row is the table's header row.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "13-gene set"
$ws.Range("C2").Value = "17-gene set"
$ws.Range("D2").Value = "30-gene set"
$ws.Range("E2").Value = "Dataset"

# ---------------------------------------------------------------------
# 2. Study / author label renames (column E), contiguous blocks.
# ---------------------------------------------------------------------
$ws.Range("E3:E41").Value = "OGR25-BTB"
$ws.Range("E42:E46").Value = "MCL14-BTB"
$ws.Range("E47:E52").Value = "WIA20-BTB"
$ws.Range("E53:E70").Value = "MCL21-BTB"

# ---------------------------------------------------------------------
# 3. Updated 13/17/30-gene set classification values (columns B,C,D),
#    rows 3-70.
# ---------------------------------------------------------------------
$bVals = @(1,0,1,0,0,1,0,0,0,1,1,0,0,1,1,0,1,1,1,0,1,0,1,1,1,1,1,1,0,1,1,1,1,1,1,1,0,0,0,0,0,0,1,1,0,0,0,1,0,1,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$cVals = @(0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,1,0,0,0,0,0,1,1,1,1,0,1,0,1,1,1,1,1,1,1,0,0,1,0,0,0,1,1,0,1,0,1,0,1,0,0,0,1,1,0,1,1,1,1,1,1,1,1,1,1,1,1)
$dVals = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,0,1,0,1,1,1,1,1,1,1,0,0,0,0,0,0,1,1,0,1,0,1,0,1,1,0,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)

for ($i = 0; $i -lt $bVals.Length; $i++) {
    $r = $i + 3
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
    $ws.Cells.Item($r, 3).Value = $cVals[$i]
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
}

# ---------------------------------------------------------------------
# 4. Row heights: row 2 (header) and rows 3-70 (data) grow to 24.95pt.
# ---------------------------------------------------------------------
$ws.Range("A2:A70").RowHeight = 24.95

# ---------------------------------------------------------------------
# 5. Drop the explicit cell style override on B/D (and B/C) so the
#    classification columns fall back to the workbook default style.
# ---------------------------------------------------------------------
$ws.Range("B2:D70").Style = "Normal"

# ---------------------------------------------------------------------
# 6. Clear the stray percentage value that used to sit in O68.
# ---------------------------------------------------------------------
$ws.Range("O68").ClearContents()

# ---------------------------------------------------------------------
# 7. Selection moves to F22.
# ---------------------------------------------------------------------
$ws.Range("F22").Select()
